$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the records (entire rows) whose Country (column A) is
# "Bosnia and Herzegovina" or "Timor-Leste" -- population data for these
# two countries was not recorded.
$countriesToRemove = @("Bosnia and Herzegovina", "Timor-Leste")

foreach ($country in $countriesToRemove) {
    $found = $ws.Columns.Item(1).Find($country)
    if ($found -ne $null) {
        $ws.Rows.Item($found.Row).Delete()
    }
}

$ws.Range("F11").Select()
